$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark from its current location (right before
#    the "Basically," run, inside the "//Hotfix" paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. The final paragraph currently reads "Pending"; it should become
#    "Ongoing" and gain the (collapsed) "_GoBack" bookmark immediately after
#    the new text, before the paragraph mark.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.MoveEnd(1, -1)      # exclude the trailing paragraph mark
if ($lastRange.Text -eq "Pending") {
    $lastRange.Text = "Ongoing"
}

# A collapsed (zero-length) Bookmarks.Add placed exactly on a paragraph
# boundary misbehaves, so anchor the bookmark using a temporary character,
# then remove the character again once the bookmark is in place.
$endRange = $d.Range($lastRange.End, $lastRange.End)
$endRange.InsertAfter("Z")

$bookmarkRange = $d.Range($endRange.Start, $endRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$bm = $d.Bookmarks("_GoBack")
$placeholder = $d.Range($bm.End, $bm.End + 1)
$placeholder.Delete()
